$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows for removed worker periods (rows 17-21)
$ws.Range("B17:J21").EntireRow.Delete() | Out-Null

$ws.Range("G16").Value = 1423500
$ws.Range("E11").Value = 40000
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
